$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing row (25) down into the new row (26)
# so the new row matches the visual style used by the rest of the table.
$ws.Range("A25:C25").Copy() | Out-Null
$ws.Range("A26:C26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row: Question No, GFG/LC, Question
$ws.Cells.Item(26, 1).Value = 226
$ws.Cells.Item(26, 2).Value = "GFG"
$ws.Cells.Item(26, 3).Value = "Invert Binary Tree(Mirror Tree)"

# Update the active selection to reflect the new end-of-table cell
$ws.Range("C26").Select()
